# Applies a cyclic swap of artfynd record data among rows 9-12 and rows 23-30
# (the row *numbers*, styles and structure stay put; only the per-record
# field values - A,B,D,E,F,G,H,Q,R - move between rows), matching the
# upstream re-export that reshuffled which physical row holds which record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 9-12
# Row 9 <- data previously on row 11
$ws.Range("A9").Value = 112044178
$ws.Range("B9").Value = 96348
$ws.Range("D9").Value = 'VU'
$ws.Range("E9").Value = 220787
$ws.Range("F9").Value = 'Knärot'
$ws.Range("G9").Value = 'Goodyera repens'
$ws.Range("H9").Value = '(L.) R. Br.'
$ws.Range("Q9").Value = 554737.3165861247
$ws.Range("R9").Value = 6697620.516129929

# Row 10 <- data previously on row 12
$ws.Range("A10").Value = 112044172
$ws.Range("B10").Value = 90658
$ws.Range("D10").Value = 'NT'
$ws.Range("E10").Value = 4361
$ws.Range("F10").Value = 'Orange taggsvamp'
$ws.Range("G10").Value = 'Hydnellum aurantiacum'
$ws.Range("H10").Value = '(Batsch:Fr.) P.Karst.'
$ws.Range("Q10").Value = 554721.7291097966
$ws.Range("R10").Value = 6697603.976012163

# Row 11 <- data previously on row 9
$ws.Range("A11").Value = 112044164
$ws.Range("B11").Value = 88924
$ws.Range("D11").Value = 'LC'
$ws.Range("E11").Value = 256703
$ws.Range("F11").Value = 'Tallfingersvamp'
$ws.Range("G11").Value = 'Ramaria eosanguinea'
$ws.Range("H11").Value = 'R.H.Petersen'
$ws.Range("Q11").Value = 554724.8915453397
$ws.Range("R11").Value = 6697591.177257041

# Row 12 <- data previously on row 10
$ws.Range("A12").Value = 112044176
$ws.Range("B12").Value = 96348
$ws.Range("D12").Value = 'VU'
$ws.Range("E12").Value = 220787
$ws.Range("F12").Value = 'Knärot'
$ws.Range("G12").Value = 'Goodyera repens'
$ws.Range("H12").Value = '(L.) R. Br.'
$ws.Range("Q12").Value = 554725.1969658234
$ws.Range("R12").Value = 6697570.923917417

# Rows 23-30
# Row 23 <- data previously on row 27
$ws.Range("A23").Value = 112044169
$ws.Range("B23").Value = 89845
$ws.Range("D23").Value = 'VU'
$ws.Range("E23").Value = 1209
$ws.Range("F23").Value = 'Rynkskinn'
$ws.Range("G23").Value = 'Phlebia centrifuga'
$ws.Range("H23").Value = 'P.Karst.'
$ws.Range("Q23").Value = 554764.5822306949
$ws.Range("R23").Value = 6697617.468857886

# Row 24 <- data previously on row 28
$ws.Range("A24").Value = 112044155
$ws.Range("B24").Value = 89405
$ws.Range("D24").Value = 'NT'
$ws.Range("E24").Value = 1202
$ws.Range("F24").Value = 'Ullticka'
$ws.Range("G24").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H24").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q24").Value = 554761.4420383665
$ws.Range("R24").Value = 6697628.785556345

# Row 25 <- data previously on row 30
$ws.Range("A25").Value = 112044157
$ws.Range("B25").Value = 89405
$ws.Range("D25").Value = 'NT'
$ws.Range("E25").Value = 1202
$ws.Range("F25").Value = 'Ullticka'
$ws.Range("G25").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H25").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q25").Value = 554763.6075584656
$ws.Range("R25").Value = 6697616.465966055

# Row 26 <- data previously on row 24
$ws.Range("A26").Value = 112044179
$ws.Range("B26").Value = 96348
$ws.Range("D26").Value = 'VU'
$ws.Range("E26").Value = 220787
$ws.Range("F26").Value = 'Knärot'
$ws.Range("G26").Value = 'Goodyera repens'
$ws.Range("H26").Value = '(L.) R. Br.'
$ws.Range("Q26").Value = 554794.5978541592
$ws.Range("R26").Value = 6697596.181857388

# Row 27 <- data previously on row 23
$ws.Range("A27").Value = 112044180
$ws.Range("B27").Value = 96348
$ws.Range("D27").Value = 'VU'
$ws.Range("E27").Value = 220787
$ws.Range("F27").Value = 'Knärot'
$ws.Range("G27").Value = 'Goodyera repens'
$ws.Range("H27").Value = '(L.) R. Br.'
$ws.Range("Q27").Value = 554838.8808180906
$ws.Range("R27").Value = 6697580.545608173

# Row 28 <- data previously on row 29
$ws.Range("A28").Value = 112044154
$ws.Range("B28").Value = 89405
$ws.Range("D28").Value = 'NT'
$ws.Range("E28").Value = 1202
$ws.Range("F28").Value = 'Ullticka'
$ws.Range("G28").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H28").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q28").Value = 554768.2498027334
$ws.Range("R28").Value = 6697636.793724483

# Row 29 <- data previously on row 25
$ws.Range("A29").Value = 112044168
$ws.Range("B29").Value = 89845
$ws.Range("D29").Value = 'VU'
$ws.Range("E29").Value = 1209
$ws.Range("F29").Value = 'Rynkskinn'
$ws.Range("G29").Value = 'Phlebia centrifuga'
$ws.Range("H29").Value = 'P.Karst.'
$ws.Range("Q29").Value = 554760.6686302377
$ws.Range("R29").Value = 6697614.44524945

# Row 30 <- data previously on row 26
$ws.Range("A30").Value = 112044156
$ws.Range("B30").Value = 89405
$ws.Range("D30").Value = 'NT'
$ws.Range("E30").Value = 1202
$ws.Range("F30").Value = 'Ullticka'
$ws.Range("G30").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H30").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q30").Value = 554761.6582123891
$ws.Range("R30").Value = 6697614.460182385

# The stray placeholder cell AF9 (empty "Bestämningsmetod" value) travelled
# with its row's data to row 11; row 9 now has no such placeholder.
$ws.Range("AF9").ClearContents()

